$d = $word.ActiveDocument

# --- 1) OLEObject ObjectID: find the paragraph that owns the EMBED field ---
$oleIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $oleIdx = $i
    }
}
if ($oleIdx -eq -1) { throw "OLE paragraph not found" }
$oleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w14:paraId="0406775E" w14:textId="77777777" w:rsidR="006E3F62" w:rsidRDefault="002B3919"><w:pPr><w:snapToGrid w:val="0"/><w:spacing w:before="960" w:line="280" w:lineRule="exact"/><w:jc w:val="both"/></w:pPr><w:r><w:object w:dxaOrig="683" w:dyaOrig="671" w14:anchorId="1EA521A4"><v:shape id="ole_rId4" o:spid="_x0000_i1026" style="width:60pt;height:59.25pt" coordsize="" o:spt="100" adj="0,,0" path="" stroked="f"><v:stroke joinstyle="miter"/><v:imagedata r:id="rId9" o:title=""/><v:formulas/><v:path o:connecttype="segments"/></v:shape><o:OLEObject Type="Embed" ProgID="PBrush" ShapeID="ole_rId4" DrawAspect="Content" ObjectID="_1667310275" r:id="rId10"/></w:object></w:r></w:p>'
$d.Paragraphs.Item($oleIdx).Range.InsertXML($oleXml)

# --- 2) OVR paragraph: locate by unique text marker ---
$ovrIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("OVR - ")) {
        $ovrIdx = $i
    }
}
if ($ovrIdx -eq -1) { throw "OVR paragraph not found" }
$ovrXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="789DF65F" w14:textId="28B086B5" w:rsidR="006E3F62" w:rsidRDefault="002B3919"><w:pPr><w:pStyle w:val="Legenda"/></w:pPr><w:r><w:t xml:space="preserve">        OVR - DOSSIÊ nº  {</w:t></w:r><w:r w:rsidR="00B831DA"><w:t>ovr_</w:t></w:r><w:r><w:t>id}</w:t></w:r></w:p>'
$d.Paragraphs.Item($ovrIdx).Range.InsertXML($ovrXml)

# --- 3) {responsavel} paragraph ---
$respIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("{responsavel}")) {
        $respIdx = $i
    }
}
if ($respIdx -eq -1) { throw "responsavel paragraph not found" }
$respXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7147259C" w14:textId="77777777" w:rsidR="006E3F62" w:rsidRDefault="002B3919"><w:pPr><w:snapToGrid w:val="0"/><w:spacing w:before="57"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>auditor_responsavel</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p>'
$d.Paragraphs.Item($respIdx).Range.InsertXML($respXml)

# --- 4) {user_name} paragraph ---
$userIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("{user_name}")) {
        $userIdx = $i
    }
}
if ($userIdx -eq -1) { throw "user_name paragraph not found" }
$userXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3B0690FA" w14:textId="77777777" w:rsidR="006E3F62" w:rsidRDefault="002B3919"><w:pPr><w:snapToGrid w:val="0"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>responsavel</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/></w:rPr><w:t>}</w:t></w:r></w:p>'
$d.Paragraphs.Item($userIdx).Range.InsertXML($userXml)

Write-Output ("oleIdx=" + $oleIdx + " ovrIdx=" + $ovrIdx + " respIdx=" + $respIdx + " userIdx=" + $userIdx)
